# Adds the "2023" column (S) to the data table, mirroring the formatting
# of the existing "2022" column (R), and adjusts the header merge,
# selection and window-size bookkeeping to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the header merge A1:R1 -> A1:S1 -----------------------------
$ws.Range("A1:R1").UnMerge()
$ws.Range("A1:S1").Merge()

# --- 2. Copy formatting from column R (rows 1-15) into column S ------------
$ws.Range("R1:R15").Copy()
$ws.Range("S1:S15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2b. Copy formatting from column T (rows 2-15) into column U -----------
# (row 1 never had a T1/U1 cell - the header merge only spans A1:S1)
$ws.Range("T2:T15").Copy()
$ws.Range("U2:U15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Fill in the 2023 year header and data values ------------------------
$ws.Range("S3").Value = 2023

$ws.Range("S4").Value = 1124.4000000000001
$ws.Range("S5").Value = 170.8
$ws.Range("S6").Value = 7146
$ws.Range("S7").Value = 4928
$ws.Range("S8").Value = 650.20000000000005
$ws.Range("S9").Value = 35
$ws.Range("S10").Value = 38.9
$ws.Range("S11").Value = 135.69999999999999
$ws.Range("S12").Value = 10.7
$ws.Range("S13").Value = 1014
$ws.Range("S14").Value = 981.2

# --- 4. Update sheet view: scroll/selection ---------------------------------
$ws.Range("B1").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("S3:S14").Select() | Out-Null
$excel.ActiveCell = $ws.Range("S3")

# --- 5. Window size bookkeeping ---------------------------------------------
$excel.ActiveWindow.Width = 741
$excel.ActiveWindow.Height = 309
